# DRomicsTODOlist.xlsx - "Add of the function trendplot()" edit
#
# The TODO item "58. ajouter une fonction trendplot (cf. article diuron)"
# (together with its owner "ML") is moved from the "short term" sheet
# (sheet2) to the "done" sheet (sheet1), i.e. the task has been completed.

$wb = $excel.ActiveWorkbook

$doneSheet  = $wb.Worksheets.Item(1)   # "done"
$shortSheet = $wb.Worksheets.Item(2)   # "short term"

# ---------------------------------------------------------------------
# 1) Append the finished task as a new last row (row 37) on the "done"
#    sheet, re-using the same formatting as the previous last row (36).
# ---------------------------------------------------------------------
$doneSheet.Range("A36").Copy()
$doneSheet.Range("A37").PasteSpecial(-4122)   # xlPasteFormats
$doneSheet.Range("A37").Value = "58. ajouter une fonction trendplot (cf. article diuron)"
$doneSheet.Range("B37").Value = "ML"

# ---------------------------------------------------------------------
# 2) Remove that same task from the "short term" sheet: it used to live
#    on row 21 (with an empty formatting-only row 22 right below it).
# ---------------------------------------------------------------------
$shortSheet.Rows.Item(21).Delete()
$shortSheet.Rows.Item(21).Delete()

# ---------------------------------------------------------------------
# 3) Row 13 on the "short term" sheet switches from the "to do" style to
#    the "in progress" style (same one used by rows 16/17/20), keeping
#    its existing text/value untouched.
# ---------------------------------------------------------------------
$shortSheet.Range("A16").Copy()
$shortSheet.Range("A13").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 4) Update the views: "short term" becomes the active/visible tab, with
#    a fresh, simple single-cell selection, while "done" keeps a
#    multi-row selection around the newly added row but is no longer the
#    selected tab.
# ---------------------------------------------------------------------
$doneSheet.Activate()
$doneSheet.Range("A43:A46").Select()

$shortSheet.Activate()
$shortSheet.Range("A13").Select()
